# Update database and roll the reporting window forward by one fiscal year,
# and fix the "read_price" (D15) cell to carry a real numeric 0 instead of a
# placeholder "-" string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------------------
# Row 8: fiscal-period column headers (shift one year forward, add 1401/12)
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---------------------------------------------------------------------------
# Row 9: publish-date column headers (shift forward, add latest release date)
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "1399-03-19 (9)"
$ws.Range("E9").Value = "1400-02-21 (8)"
$ws.Range("F9").Value = "1401-03-04 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "1402-02-28"

# ---------------------------------------------------------------------------
# Rows 11-27: financial data, shifted one column to the left (dropping the
# oldest fiscal year) with a freshly computed rightmost (H) column.
# Each inner array is: row, D, E, F, G, H
# ---------------------------------------------------------------------------
$cols = @("D", "E", "F", "G", "H")

$data = @(
    @(11, 4892461, 7408342, 8800845, 14540574, 23131160),
    @(12, -2663735, -3795904, -4503737, -7957214, -13518131),
    @(13, 2228726, 3612438, 4297108, 6583360, 9613029),
    @(14, -198540, -133977, -140369, -305190, -339554),
    @(15, 0, 0, 0, 0, 0),
    @(16, 243303, -25189, 91002, -83191, -154184),
    @(17, 2273489, 3453272, 4247741, 6194979, 9119291),
    @(18, -435598, -686481, -865631, -1545549, -1937565),
    @(19, 31788, 390983, 742406, 659657, 837932),
    @(20, 1869679, 3157774, 4124516, 5309087, 8019658),
    @(21, -356678, -663097, -731176, -808726, -899914),
    @(22, 1513001, 2494677, 3393340, 4500361, 7119744),
    @(23, 0, 0, 0, 0, 0),
    @(24, 1513001, 2494677, 3393340, 4500361, 7119744),
    @(25, 813, 1341, 1542, 1461, 1082),
    @(26, 1860000, 1860000, 2200000, 3080000, 6580000),
    @(27, 230, 379, 516, 684, 1082)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $row).Value = $entry[$i + 1]
    }
}
